# Error Calculations and Plots
# Apply the edits described by the diff to missing_data.xlsx (Sheet1).
#
# The diff removes two entire data rows ("RM 232" and "SC 92") -- everything
# below each shifts up -- and then tweaks a handful of individual cells
# (clearing some, filling in previously-missing values in others).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two removed rows (delete the lower row first so the
#        higher row's index doesn't move before we get to it). ---
# Row 28 = "SC 92" (before any deletions)
$ws.Rows(28).Delete()
# Row 26 = "RM 232" (before any deletions; still 26 since it is above row 28)
$ws.Rows(26).Delete()

# --- 2. Apply the remaining individual cell edits (row numbers below are
#        the numbers AFTER the two rows above have been removed). ---

# Row 5 (RM 14): column D (header "C") is now missing
$ws.Range("D5").ClearContents()

# Row 11 (RM 58): column D (header "C") is now filled in
$ws.Range("D11").Value = -15.5

# Row 12 (RM 81): column F is now missing
$ws.Range("F12").ClearContents()

# Row 13 (RM 88): column F is now filled in
$ws.Range("F13").Value = 17.1

# Row 16 (RM 103): column F is now filled in
$ws.Range("F16").Value = 17.34

# Row 17 (RM 116): column F is now filled in
$ws.Range("F17").Value = 17.78

# Row 19 (RM 125): column C (header "B") is now filled in, column D (header "C") is now missing
$ws.Range("C19").Value = 13.2
$ws.Range("D19").ClearContents()

# Row 20 (RM 134): column F is now filled in
$ws.Range("F20").Value = 17.73

# Row 21 (RM 135): column C (header "B") is now missing
$ws.Range("C21").ClearContents()

# Row 22 (RM 138): column F is now filled in
$ws.Range("F22").Value = 16.81

# Row 23 (RM 140): column C (header "B") is now filled in
$ws.Range("C23").Value = 12.2

# Row 24 (RM 142a): column F is now missing
$ws.Range("F24").ClearContents()

# Row 25 (RM 145): column D (header "C") is now filled in, column F is now missing
$ws.Range("D25").Value = -15.5
$ws.Range("F25").ClearContents()

# Row 27 (SC 101): column C (header "B") is now missing
$ws.Range("C27").ClearContents()

# Row 28 (SC 105): column E (header "D") is now filled in, column F is now missing
$ws.Range("E28").Value = -5.9
$ws.Range("F28").ClearContents()

# Row 29 (SC 119): column D (header "C") is now missing, column F is now missing
$ws.Range("D29").ClearContents()
$ws.Range("F29").ClearContents()

# Row 30 (SC 120): column F is now missing
$ws.Range("F30").ClearContents()

# Row 32 (SC 193): column E (header "D") is now missing
$ws.Range("E32").ClearContents()

# Row 33 (SC 232): column C (header "B") is now filled in
$ws.Range("C33").Value = 10.4
